# Update gh-pages output ("want to go" counts refreshed + a new exhibition
# row added) on the two sheets that list every exhibition: "展览" and the
# roll-up "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the last populated row (column A holds a sequential index).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    # ---- refresh the "想去人数" (want-to-go count) column F ----------------
    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        switch ($name) {
            "上饶·囧喵喵国风动漫展"                                  { $ws.Cells.Item($r, 6).Value = 444 }
            "南昌·第一届异次元动漫嘉年华"                             { $ws.Cells.Item($r, 6).Value = 1832 }
            "赣州·第二届异次元动漫嘉年华"                             { $ws.Cells.Item($r, 6).Value = 1446 }
            "信丰·七夕节UPUP动漫展"                                  { $ws.Cells.Item($r, 6).Value = 135 }
            "南昌·W·MEETING动漫游戏盛典"                             { $ws.Cells.Item($r, 6).Value = 1716 }
            "南昌·花绒万兽第二聚"                                    { $ws.Cells.Item($r, 6).Value = 139 }
            "吉安·WF无线次元新星动漫博览会"                           { $ws.Cells.Item($r, 6).Value = 648 }
            "赣州·十万伏特-星铁&音乐 次元音乐同人only2.0"              { $ws.Cells.Item($r, 6).Value = 56 }
            "抚州·逆光ZERO动漫游戏展"                                { $ws.Cells.Item($r, 6).Value = 75 }
            "萍乡·夏花Flower·2024夏季国漫展"                         { $ws.Cells.Item($r, 6).Value = 138 }
            "上饶·次元重现夏日嘉年华（取消）"                          { $ws.Cells.Item($r, 6).Value = 121 }
            "乐平·CY境界次元第三届动漫游戏庆典"                       { $ws.Cells.Item($r, 6).Value = 63 }
            "南昌·CM03·配音演员孙路路专场见面会"                      { $ws.Cells.Item($r, 6).Value = 100 }
            "南昌·CM03动漫游戏博览会"                                { $ws.Cells.Item($r, 6).Value = 4480 }
            "九江·如梦令国潮动漫节"                                  { $ws.Cells.Item($r, 6).Value = 35 }
            "南昌·第四届龙年动漫展——暑假最后的狂欢"                    { $ws.Cells.Item($r, 6).Value = 802 }
            "赣州·第五人格only"                                     { $ws.Cells.Item($r, 6).Value = 96 }
            "南昌·Sunflower Garden动漫游戏展"                       { $ws.Cells.Item($r, 6).Value = 2141 }
            "南昌·第一届哥布林动漫游戏展——开学季&贺中秋"                { $ws.Cells.Item($r, 6).Value = 69 }
            "南昌·萌卡动漫展"                                       { $ws.Cells.Item($r, 6).Value = 2012 }
        }
    }

    # ---- insert the new "南昌·Aud中秋动漫嘉年华" row just above the final row --
    $insertAt = $lastRow

    # Insert a blank row first, then clone column-A's bordered/bold/centered
    # look from the row above onto the new blank index cell.
    $ws.Rows.Item($insertAt).Insert()
    $ws.Range("A" + ($insertAt - 1)).Copy()
    $ws.Range("A" + $insertAt).PasteSpecial(-4122)

    $ws.Cells.Item($insertAt, 1).Value = $insertAt - 1

    # Column B holds a plain-text date string ("2024-09-17"), not a real
    # date value, so force text formatting before assigning it, then strip
    # the formatting override back off by copying the "no override" look
    # from an existing date cell in the column.
    $dateCell = $ws.Cells.Item($insertAt, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2024-09-17"
    $ws.Cells.Item(2, 2).Copy()
    $dateCell.PasteSpecial(-4122)

    $ws.Cells.Item($insertAt, 3).Value = "南昌·Aud中秋动漫嘉年华"
    $ws.Cells.Item($insertAt, 4).Value = "青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK"
    $ws.Cells.Item($insertAt, 5).Value = "2024.09.17 10:00-09.17 17:00"
    $ws.Cells.Item($insertAt, 6).Value = 2
    $ws.Cells.Item($insertAt, 7).Value = 29.9
    $ws.Cells.Item($insertAt, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90329"
    $ws.Cells.Item($insertAt, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/pbU7Eftp1722660514298.jpeg"

    # Fix the index number of the row that just got pushed down.
    $ws.Cells.Item($insertAt + 1, 1).Value = $insertAt
}
